$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.919.67'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.629.78'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.83'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.36'
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.256'
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.861.05'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.627.70'
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('E14').Value = '  -1.43%  '
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.913.06'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.29'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.64'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.25'
$ws.Range('E23').Value = '  -4.42%  '
$ws.Range('E24').Value = '  -1.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.82'
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.92'
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('E32').Value = '  +2.35%  '
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.402.09'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('E36').Value = '  +10.93%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('E38').Value = '  +1.97%  '
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('E43').Value = '  +0.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.26'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.48'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('E46').Value = '  -0.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.770.75'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.11'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('E49').Value = '  -2.48%  '
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('E51').Value = '  -0.32%  '
